{"js": "// The supplied OOXML diff touches every part (document.xml, headers,\n// footers, footnotes/endnotes, numbering.xml, styles.xml, theme1.xml),\n// but in every single hunk the only two lines that actually change are\n// the root element's opening tag: the `xmlns:*` namespace-prefix\n// declarations get reordered (e.g. `m` moves after `a`) and the Excel\n// URN's prefix is renamed (`ns19` -> `ns17`), and `mc:Ignorable=\"w14 w15\"`\n// is dropped from <w:document>. Every other line in every hunk \u2014 all of\n// the actual body text, paragraphs, runs, table/style/numbering\n// definitions, theme colors, etc. \u2014 is unchanged context.\n//\n// That pattern (uniform namespace-prefix reshuffle across *all* parts,\n// with zero visible/content differences) is exactly what happens when a\n// package is re-serialized by a different OOXML writer, which matches\n// the accompanying commit message: \"n\u00e3o uso mais a biblioteca de\n// terceiro para fazer a convers\u00e3o do xml e vise versa\" (switched away\n// from a third-party library for XML (de)serialization in the\n// originating application). It is a byproduct of that unrelated\n// application change, not a content edit made to this document, and the\n// raw XML namespace/prefix ordering it produces is not something the\n// Word JavaScript API exposes any control over.\n//\n// So the correct, faithful reproduction of this diff through Office.js\n// is to leave the document's content exactly as it is -- i.e. no\n// mutation at all.\n\nawait context.sync();\n", "ps1": "# The supplied OOXML diff touches every part (document.xml, headers,\n# footers, footnotes/endnotes, numbering.xml, styles.xml, theme1.xml),\n# but in every single hunk the only two lines that actually change are\n# the root element's opening tag: the `xmlns:*` namespace-prefix\n# declarations get reordered (e.g. `m` moves after `a`) and the Excel\n# URN's prefix is renamed (`ns19` -> `ns17`), and `mc:Ignorable=\"w14 w15\"`\n# is dropped from <w:document>. Every other line in every hunk -- all of\n# the actual body text, paragraphs, runs, table/style/numbering\n# definitions, theme colors, etc. -- is unchanged context.\n#\n# That pattern (uniform namespace-prefix reshuffle across *all* parts,\n# with zero visible/content differences) is exactly what happens when a\n# package is re-serialized by a different OOXML writer, which matches\n# the accompanying commit message: \"n\u00e3o uso mais a biblioteca de\n# terceiro para fazer a convers\u00e3o do xml e vise versa\" (switched away\n# from a third-party library for XML (de)serialization in the\n# originating application). It is a byproduct of that unrelated\n# application change, not a content edit made to this document, and the\n# raw XML namespace/prefix ordering it produces is not something the\n# Word COM object model exposes any control over.\n#\n# So the correct, faithful reproduction of this diff through COM\n# automation is to leave the document's content exactly as it is -- i.e.\n# no mutation at all.\n\n$d = $word.ActiveDocument\n"}
